$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 89
$ws.Range("E3").Value = 31
$ws.Range("E10").Value = 404
$ws.Range("F10").Value = 198
$ws.Range("H10").Value = 198
$ws.Range("E11").Value = 275
$ws.Range("E12").Value = 391
$ws.Range("F12").Value = 215
$ws.Range("H12").Value = 215
$ws.Range("E13").Value = 105
$ws.Range("E15").Value = 129
$ws.Range("F20").Value = 27
$ws.Range("H20").Value = 27
$ws.Range("E23").Value = 170
$ws.Range("E24").Value = 173
$ws.Range("E25").Value = 205
$ws.Range("E26").Value = 117
$ws.Range("F26").Value = 72
$ws.Range("H26").Value = 72
$ws.Range("E27").Value = 265
$ws.Range("E28").Value = 160
$ws.Range("F30").Value = 97
$ws.Range("H30").Value = 97
$ws.Range("E32").Value = 162
$ws.Range("F32").Value = 93
$ws.Range("H32").Value = 93
$ws.Range("E33").Value = 247
$ws.Range("E34").Value = 177
$ws.Range("F34").Value = 107
$ws.Range("H34").Value = 107
$ws.Range("E35").Value = 116
$ws.Range("F35").Value = 74
$ws.Range("H35").Value = 74
$ws.Range("E37").Value = 130
$ws.Range("E38").Value = 80
$ws.Range("E39").Value = 161
$ws.Range("E40").Value = 219
$ws.Range("F40").Value = 99
$ws.Range("H40").Value = 99
$ws.Range("E41").Value = 324
$ws.Range("E42").Value = 293
$ws.Range("E44").Value = 256
$ws.Range("E45").Value = 115
$ws.Range("F45").Value = 56
$ws.Range("H45").Value = 56
$ws.Range("E46").Value = 256
$ws.Range("E47").Value = 366
$ws.Range("E48").Value = 173
$ws.Range("E49").Value = 248
$ws.Range("E50").Value = 211
$ws.Range("F50").Value = 90
$ws.Range("H50").Value = 90
$ws.Range("E51").Value = 199
